$d = $word.ActiveDocument

$replacements = @(
    @{old="66×25=1650"; new="99×26=2574"},
    @{old="98×27=2646"; new="39×89=3471"},
    @{old="76×27=2052"; new="71×82=5822"},
    @{old="41×29=1189"; new="71×36=2556"},
    @{old="70×12=840";  new="27×12=324"},
    @{old="34×54=1836"; new="90×22=1980"},
    @{old="98×52=5096"; new="15×82=1230"},
    @{old="75×34=2550"; new="19×43=817"},
    @{old="44×71=3124"; new="84×62=5208"},
    @{old="80×67=5360"; new="34×20=680"},
    @{old="29×18=522";  new="22×63=1386"},
    @{old="96×62=5952"; new="41×20=820"},
    @{old="30×51=1530"; new="66×37=2442"},
    @{old="75×45=3375"; new="35×22=770"},
    @{old="66×95=6270"; new="87×60=5220"},
    @{old="57×17=969";  new="39×22=858"},
    @{old="16×91=1456"; new="44×90=3960"},
    @{old="88×81=7128"; new="53×11=583"},
    @{old="13×39=507";  new="26×64=1664"},
    @{old="68×16=1088"; new="11×57=627"},
    @{old="30×64=1920"; new="21×85=1785"},
    @{old="99×43=4257"; new="86×27=2322"},
    @{old="66×32=2112"; new="50×81=4050"},
    @{old="46×70=3220"; new="42×35=1470"},
    @{old="63×97=6111"; new="12×36=432"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
